# Edit: fix the "Jogo da velha" bullet text (drop the redundant "que")
# and relocate the "_GoBack" bookmark (Word's last-edit-position marker)
# from the old "Pagina do jogo" spot to the new edit point.

$d = $word.ActiveDocument

# --- Step 1: remove the word "que " from
#     "...computador que possui tres niveis de dificuldades"
#     so it reads "...computador possui tres niveis de dificuldades"
$r1 = $d.Content
$found1 = $r1.Find.Execute("que possui três níveis de dificuldades", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $r1.Text = "possui três níveis de dificuldades"
}

# --- Step 2: make sure the trailing period after "dificuldades" stays in
#     its own run (matches the author's original run layout). We force the
#     split with a throw-away bookmark, then immediately delete the
#     bookmark itself (the run break it leaves behind persists).
$r2 = $d.Content
$found2 = $r2.Find.Execute("dificuldades.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $splitAt = $r2.Start + ($r2.End - $r2.Start - 1)
    $splitRange = $d.Range($splitAt, $splitAt)
    $d.Bookmarks.Add("TEMP_SPLIT_MARK", $splitRange)
    $d.Bookmarks.Item("TEMP_SPLIT_MARK").Delete()
}

# --- Step 3: drop a bookmark named "_GoBack" right before "possui" (between
#     "computador " and "possui tres niveis..."). Word only ever keeps a
#     single "_GoBack" bookmark, so adding it here automatically removes it
#     from its previous location after "Pagina do jogo".
$r3 = $d.Content
$found3 = $r3.Find.Execute("possui três níveis de dificuldades", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $r3.Collapse(1)
    $d.Bookmarks.Add("_GoBack", $r3)
}
